$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A72").Value = "2023-12-07 14:19:05"
$ws.Range("B72").Value = 0.0002
